$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A170").Value = 169
$ws.Range("B170").Value = 1
$ws.Range("C170").Value = "2024-06-18 13:15:22"
$ws.Range("D170").Value = 200
$ws.Range("E170").Value = 19

$ws.Range("A171").Value = 170
$ws.Range("B171").Value = 2
$ws.Range("C171").Value = "2024-06-18 13:15:22"
$ws.Range("D171").Value = 200
$ws.Range("E171").Value = 0
